$wb = $excel.ActiveWorkbook

# --- Sheet "mkb": remove the old deposit tables in columns A:J, keep the
#     "Вклад «МЕГА Онлайн»" / "МКБ. 30 лет" / "Все включено" block in K:U ---
$wsMkb = $wb.Worksheets.Item("mkb")
$wsMkb.Range("A1:J33").ClearContents()

# --- Sheet "vtb": rotate the minimum-deposit labels in A20:A22 ---
$wsVtb = $wb.Worksheets.Item("vtb")
$a20 = $wsVtb.Range("A20").Value2
$a21 = $wsVtb.Range("A21").Value2
$a22 = $wsVtb.Range("A22").Value2
$wsVtb.Range("A20").Value = $a22
$wsVtb.Range("A21").Value = $a20
$wsVtb.Range("A22").Value = $a21

# --- Sheet "gazprom": drop the "Розница" / "Премиум" header row ---
$wsGazprom = $wb.Worksheets.Item("gazprom")
$wsGazprom.Range("A1:J1").ClearContents()

# --- Sheet "alfa": add the "Розница" / "Премиум" header row, fix typo ---
$wsAlfa = $wb.Worksheets.Item("alfa")
$wsAlfa.Range("A1").Value = "Розница"
$wsAlfa.Range("K1").Value = "Премиум"
$wsAlfa.Range("A10").Value = "С пополнением"

# --- Sheet "psb": updated rate figures ---
$wsPsb = $wb.Worksheets.Item("psb")
$wsPsb.Range("G2").Value = "7.15% "
$wsPsb.Range("E10").Value = "5% "
$wsPsb.Range("F11").Value = "5% "
$wsPsb.Range("F12").Value = "5% "
$wsPsb.Range("F13").Value = "5% "
